# LoginData.xlsx update:
#  - A2/C2 test data refreshed to a new "mka" user (replacing the old "test"/"root root" entry)
#  - D3:E5 (Assertation/Status for rows 3-5) cleared out
#  - Hyperlinks on column A rebuilt so A2 now points at the new mka999@gmail.com address
#  - Selection left on B2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the login/email test data for row 2
$ws.Range("A2").Value = "mka999@gmail.com"
$ws.Range("C2").Value = "mka mka"

# Clear the Assertation/Status columns for rows 3-5
$ws.Range("D3:E5").ClearContents()

# Rebuild the hyperlinks on column A so that A2 points to the new address
# while A3/A4/A5 keep pointing at their existing mailto addresses.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:malik999@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:Alice999@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:Bob999@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:mka999@gmail.com")

# Re-apply the built-in Hyperlink style to column A (Hyperlinks.Add leaves
# a duplicate style behind otherwise) so all four cells share the same style.
$ws.Range("A2:A5").Style = "Hyperlink"

# Leave the selection on B2, matching the saved view state.
[void]$ws.Range("B2").Select()
